# Exp4_Jenkins.docx update — resync the inline picture docPr/cNvPr "name"
# attributes (these had drifted out of sync with the actual embedded
# media part) and shrink two pictures (image8.jpg->image7.jpg and
# image1.jpg) to ~90% of their former size, matching the latest export.
#
# InlineShape has no settable .Name in the Word object model (only
# floating Shape objects expose .Name), so the rename is done by
# editing the drawing's WordOpenXML in place: capture the shape's
# Range, patch the "name=" attributes in its WordOpenXML, delete the
# range, then re-insert the patched markup at the same spot. This
# preserves the r:embed relationship (and therefore the underlying
# image part) untouched — only the cosmetic name changes.

$d = $word.ActiveDocument

function Rename-InlineShapeXml($Document, $InlineShape, $OldName, $NewName) {
    $r = $InlineShape.Range
    $xml = $r.WordOpenXML
    $oldAttr = 'name="' + $OldName + '"'
    $newAttr = 'name="' + $NewName + '"'
    $patched = $xml.Replace($oldAttr, $newAttr)
    $start = $r.Start
    $r.Delete()
    $insertAt = $Document.Range($start, $start)
    $insertAt.InsertXML($patched)
}

# Inline pictures in document order, by their current (pre-edit) docPr
# id/name, and the new name each should carry.
Rename-InlineShapeXml $d $d.InlineShapes.Item(1) "image7.png" "image8.png"
Rename-InlineShapeXml $d $d.InlineShapes.Item(2) "image2.png" "image6.png"
Rename-InlineShapeXml $d $d.InlineShapes.Item(3) "image6.png" "image2.png"
Rename-InlineShapeXml $d $d.InlineShapes.Item(4) "image9.png" "image10.png"
Rename-InlineShapeXml $d $d.InlineShapes.Item(5) "image5.png" "image3.png"
Rename-InlineShapeXml $d $d.InlineShapes.Item(7) "image4.jpg" "image9.jpg"
Rename-InlineShapeXml $d $d.InlineShapes.Item(8) "image8.jpg" "image7.jpg"
Rename-InlineShapeXml $d $d.InlineShapes.Item(10) "image3.jpg" "image5.jpg"
Rename-InlineShapeXml $d $d.InlineShapes.Item(11) "image10.jpg" "image4.jpg"

# Shrink the two re-exported pictures (now ~90.1%-90.3% of their prior
# extent) — Width/Height are in points; 1 pt = 12700 EMU.
$pic7 = $d.InlineShapes.Item(8)
$pic7.Width = 5176838 / 12700.0
$pic7.Height = 2201446 / 12700.0

$pic1 = $d.InlineShapes.Item(9)
$pic1.Width = 3439950 / 12700.0
$pic1.Height = 2058823 / 12700.0
